# Applies the "Updated items, affixes and monsters." edit to the Affixes sheet:
#  - Updates the R2:R62 "gold cost" column values to their new (much larger) amounts.
#  - Updates the sheet view's topLeftCell / selection to reflect where the author
#    left the cursor after editing (D19 top-left, Y60 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Affixes")

# New values for R2:R62, in row order.
$newValues = @(
    500, 500, 500, 500, 500, 500,
    750,
    1500, 1500, 1500, 1500,
    3000, 3000, 3000, 3000,
    5000,
    10000, 10000, 10000, 10000, 10000, 10000,
    25000, 25000, 25000, 25000,
    50000, 50000, 50000, 50000,
    100000, 100000,
    150000, 150000, 150000, 150000, 150000,
    300000, 300000, 300000, 300000,
    500000, 500000, 500000, 500000, 500000, 500000,
    750000, 750000, 750000, 750000, 750000, 750000, 750000,
    1000000,
    5000000, 5000000, 5000000, 5000000, 5000000, 5000000
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 18).Value = $newValues[$i]
}

# Update the view state (scroll position + active selection) left behind by the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("Y60").Select()
